$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value2 = 46053
$ws.Range("D8").Value2 = 157.68
$ws.Range("E8").Value2 = 149.07
$ws.Range("F8").Value2 = 159.07
$ws.Range("G8").Value2 = 149.09

$ws.Range("A9").Value2 = 46053
$ws.Range("D9").Value2 = 157.68
$ws.Range("E9").Value2 = 149.07
$ws.Range("F9").Value2 = 159.07
$ws.Range("G9").Value2 = 149.09

$ws.Range("A10").Value2 = 46053
$ws.Range("D10").Value2 = 158.41
$ws.Range("E10").Value2 = 150.67
$ws.Range("F10").Value2 = 160.67
$ws.Range("G10").Value2 = 151.07

$ws.Range("A11").Value2 = 46052
$ws.Range("D11").Value2 = 157.95
$ws.Range("E11").Value2 = 149.37
$ws.Range("F11").Value2 = 159.37
$ws.Range("G11").Value2 = 149.39

$ws.Range("A12").Value2 = 46052
$ws.Range("D12").Value2 = 157.95
$ws.Range("E12").Value2 = 149.37
$ws.Range("F12").Value2 = 159.37
$ws.Range("G12").Value2 = 149.39

$ws.Range("A13").Value2 = 46052
$ws.Range("D13").Value2 = 158.77
$ws.Range("E13").Value2 = 151.17
$ws.Range("F13").Value2 = 161.17
$ws.Range("G13").Value2 = 151.57

$ws.Range("A17").Value2 = 46053
$ws.Range("D17").Value2 = 161.98
$ws.Range("E17").Value2 = 153.33
$ws.Range("F17").Value2 = 163.33

$ws.Range("A18").Value2 = 46052
$ws.Range("D18").Value2 = 162.3
$ws.Range("E18").Value2 = 153.79
$ws.Range("F18").Value2 = 163.79

$ws.Range("A22").Value2 = 46053
$ws.Range("D22").Value2 = 158.64
$ws.Range("E22").Value2 = 150.53
$ws.Range("F22").Value2 = 160.13
$ws.Range("G22").Value2 = 151.61

$ws.Range("A23").Value2 = 46053
$ws.Range("D23").Value2 = 163.18
$ws.Range("E23").Value2 = 156.02
$ws.Range("F23").Value2 = 166.02

$ws.Range("A24").Value2 = 46053
$ws.Range("D24").Value2 = 163.34
$ws.Range("E24").Value2 = 156.67
$ws.Range("F24").Value2 = 166.67

$ws.Range("A25").Value2 = 46053
$ws.Range("D25").Value2 = 163.34
$ws.Range("E25").Value2 = 156.2
$ws.Range("F25").Value2 = 166.2
$ws.Range("G25").Value2 = 156.33

$ws.Range("A26").Value2 = 46053
$ws.Range("D26").Value2 = 162.92
$ws.Range("E26").Value2 = 157.79
$ws.Range("F26").Value2 = 167.79

$ws.Range("A27").Value2 = 46052
$ws.Range("D27").Value2 = 158.91
$ws.Range("E27").Value2 = 150.94
$ws.Range("F27").Value2 = 160.54
$ws.Range("G27").Value2 = 152.02

$ws.Range("A28").Value2 = 46052
$ws.Range("D28").Value2 = 163.54
$ws.Range("E28").Value2 = 156.51
$ws.Range("F28").Value2 = 166.51

$ws.Range("A29").Value2 = 46052
$ws.Range("D29").Value2 = 163.69
$ws.Range("E29").Value2 = 157.18
$ws.Range("F29").Value2 = 167.18

$ws.Range("A30").Value2 = 46052
$ws.Range("D30").Value2 = 163.69
$ws.Range("E30").Value2 = 156.71
$ws.Range("F30").Value2 = 166.71
$ws.Range("G30").Value2 = 156.84

$ws.Range("A31").Value2 = 46052
$ws.Range("D31").Value2 = 163.28
$ws.Range("E31").Value2 = 158.31
$ws.Range("F31").Value2 = 168.31

$ws.Range("A35").Value2 = 46053
$ws.Range("D35").Value2 = 157.34
$ws.Range("E35").Value2 = 147.7
$ws.Range("F35").Value2 = 156.7

$ws.Range("A36").Value2 = 46052
$ws.Range("D36").Value2 = 157.69
$ws.Range("E36").Value2 = 148.19
$ws.Range("F36").Value2 = 157.19

$ws.Range("A40").Value2 = 46053
$ws.Range("D40").Value2 = 162.74
$ws.Range("E40").Value2 = 154.95
$ws.Range("F40").Value2 = 164.95

$ws.Range("A41").Value2 = 46053
$ws.Range("D41").Value2 = 162.46
$ws.Range("E41").Value2 = 155.37
$ws.Range("F41").Value2 = 165.37

$ws.Range("A42").Value2 = 46052
$ws.Range("D42").Value2 = 163.68
$ws.Range("E42").Value2 = 156.15
$ws.Range("F42").Value2 = 166.15

$ws.Range("A43").Value2 = 46052
$ws.Range("D43").Value2 = 163.4
$ws.Range("E43").Value2 = 156.57
$ws.Range("F43").Value2 = 166.57

$ws.Range("A47").Value2 = 46053
$ws.Range("D47").Value2 = 157.38
$ws.Range("E47").Value2 = 149.85
$ws.Range("F47").Value2 = 159.85

$ws.Range("A48").Value2 = 46053
$ws.Range("D48").Value2 = 156.99
$ws.Range("E48").Value2 = 149.78
$ws.Range("F48").Value2 = 159.78

$ws.Range("A49").Value2 = 46052
$ws.Range("D49").Value2 = 157.57
$ws.Range("E49").Value2 = 150.14
$ws.Range("F49").Value2 = 160.14

$ws.Range("A50").Value2 = 46052
$ws.Range("D50").Value2 = 157.17
$ws.Range("E50").Value2 = 150.06
$ws.Range("F50").Value2 = 160.06

$ws.Range("A54").Value2 = 46053
$ws.Range("D54").Value2 = 171.73
$ws.Range("E54").Value2 = 163.55
$ws.Range("F54").Value2 = 173.55

$ws.Range("A55").Value2 = 46053
$ws.Range("D55").Value2 = 164.34
$ws.Range("E55").Value2 = 162.02
$ws.Range("F55").Value2 = 172.02

$ws.Range("A56").Value2 = 46053
$ws.Range("D56").Value2 = 161.36

$ws.Range("A57").Value2 = 46053
$ws.Range("D57").Value2 = 161.9
$ws.Range("E57").Value2 = 156.44

$ws.Range("A58").Value2 = 46053
$ws.Range("D58").Value2 = 157.67
$ws.Range("E58").Value2 = 152.34
$ws.Range("F58").Value2 = 162.34

$ws.Range("A59").Value2 = 46053
$ws.Range("D59").Value2 = 164.35
$ws.Range("E59").Value2 = 161.7

$ws.Range("A60").Value2 = 46052
$ws.Range("D60").Value2 = 172.13
$ws.Range("E60").Value2 = 164.15
$ws.Range("F60").Value2 = 174.15

$ws.Range("A61").Value2 = 46052
$ws.Range("D61").Value2 = 164.78
$ws.Range("E61").Value2 = 162.48
$ws.Range("F61").Value2 = 172.48

$ws.Range("A62").Value2 = 46052
$ws.Range("D62").Value2 = 161.72

$ws.Range("A63").Value2 = 46052
$ws.Range("D63").Value2 = 162.21
$ws.Range("E63").Value2 = 156.9

$ws.Range("A64").Value2 = 46052
$ws.Range("D64").Value2 = 157.98
$ws.Range("E64").Value2 = 152.8
$ws.Range("F64").Value2 = 162.8

$ws.Range("A65").Value2 = 46052
$ws.Range("D65").Value2 = 164.65
$ws.Range("E65").Value2 = 162.27
